$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Notes")

# -- Update the "Notes" worksheet table --
# Insert a new blank row just under the header row, shifting the
# existing data rows down by one, then grow the table (ListObject) to
# cover the new row.
$ws2.Range("A5:C5").Insert(-4121)  # xlShiftDown

$lo = $ws2.ListObjects.Item(1)
$lo.Resize($ws2.Range("A4:C13"))

# Populate the newly inserted row with the NMI note.
# Leading apostrophe forces the "-" to be stored as literal text
# (quote-prefixed), matching the style of the other "-" rows.
$ws2.Cells.Item(5, 1).Value = "'-"
$ws2.Cells.Item(5, 2).Value = "NMI"
$ws2.Cells.Item(5, 3).Value = "NMI is not enabled"

# The row that used to hold the "R" note (now pushed down to row 7)
# gets its Instruction column relabelled to "IFF1/2".
$ws2.Cells.Item(7, 2).Value = "IFF1/2"

# -- Switch the active sheet to "Notes" and update the selection --
$ws2.Activate() | Out-Null
$ws2.Range("B7").Select() | Out-Null
